$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Invalid (G) and Absent (H) become 1
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Row 4: Total Attendance Count (D) and Real (E) become 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1

# Row 5: Total Attendance Count (D) and Real (E) become 1
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

# Row 6: Absent (H) becomes 1
$ws.Range("H6").Value = 1

# Row 7: Absent (H) becomes 1
$ws.Range("H7").Value = 1

# Row 8: Absent (H) becomes 1
$ws.Range("H8").Value = 1

# Row 9: Total Attendance Count (D) and Real (E) become 1
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 1

# Row 10: Absent (H) becomes 1
$ws.Range("H10").Value = 1

# Row 11: Absent (H) becomes 1
$ws.Range("H11").Value = 1

# Row 12: Absent (H) becomes 1
$ws.Range("H12").Value = 1

# Row 13: Absent (H) becomes 1
$ws.Range("H13").Value = 1

# Row 14: Absent (H) becomes 1
$ws.Range("H14").Value = 1

# Row 15: Absent (H) becomes 1
$ws.Range("H15").Value = 1

# Row 16: Absent (H) becomes 1
$ws.Range("H16").Value = 1

# Row 17: Absent (H) becomes 1
$ws.Range("H17").Value = 1

# Row 18: Absent (H) becomes 1
$ws.Range("H18").Value = 1
